$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data between row 2 and row 3 ---
# Columns kept in place: A (seq id), C (Div), D (Div Original Name), E (Date), F (HomeTeam), G (AwayTeam)
# Columns swapped: B (id) and H..AC (match stats)
$swapCols = @(2) + (8..29)

$row2vals = @{}
$row3vals = @{}
foreach ($c in $swapCols) {
    $row2vals[$c] = $ws.Cells.Item(2, $c).Value()
    $row3vals[$c] = $ws.Cells.Item(3, $c).Value()
}
foreach ($c in $swapCols) {
    $ws.Cells.Item(2, $c).Value = $row3vals[$c]
    $ws.Cells.Item(3, $c).Value = $row2vals[$c]
}

# --- Rotate match data across rows 11, 12, 13 ---
# row11 <- old row13 ; row12 <- old row11 ; row13 <- old row12
# Columns kept in place: A, C, D, E, F (HomeTeam)
# Columns rotated: B (id), G (AwayTeam) and H..AC (match stats)
$rotCols = @(2, 7) + (8..29)

$row11vals = @{}
$row12vals = @{}
$row13vals = @{}
foreach ($c in $rotCols) {
    $row11vals[$c] = $ws.Cells.Item(11, $c).Value()
    $row12vals[$c] = $ws.Cells.Item(12, $c).Value()
    $row13vals[$c] = $ws.Cells.Item(13, $c).Value()
}
foreach ($c in $rotCols) {
    $ws.Cells.Item(11, $c).Value = $row13vals[$c]
    $ws.Cells.Item(12, $c).Value = $row11vals[$c]
    $ws.Cells.Item(13, $c).Value = $row12vals[$c]
}
